$d = $word.ActiveDocument

# The document ends with a centered paragraph:
#   "================="  (bold/underlined run)
#   <manual page break>  (bold/underlined run containing only <w:br w:type="page"/>)
#   " "                  (plain run, carries a stale <w:lastRenderedPageBreak/>)
#
# The edit removes the manual page break entirely and drops the now-stale
# "last rendered page break" marker from the trailing space run, while
# leaving everything else (including the space itself) untouched.

$markerRng = $d.Content
$found = $markerRng.Find.Execute("=================", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)

if ($found) {
    $afterMarker = $markerRng.End

    # Delete the manual page break character that immediately follows the
    # marker text; this removes the whole run that contains <w:br w:type="page"/>.
    $breakRng = $d.Range($afterMarker, $afterMarker + 1)
    $breakRng.Delete()

    # The trailing single-space run now sits right after the marker run.
    # It still carries the stale <w:lastRenderedPageBreak/> flag left over
    # from when the manual page break used to precede it. Force the run to
    # be rebuilt (change its text, then change it back) so the flag is
    # dropped; using the Text property this way regenerates the run without
    # merging it into the differently-formatted marker run.
    $spaceRng = $d.Range($afterMarker, $afterMarker + 1)
    $spaceRng.Text = "X"
    $spaceRng2 = $d.Range($afterMarker, $afterMarker + 1)
    $spaceRng2.Text = " "
}
